$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 386.0435
$ws.Range("J2").Value = 622.1429000000001
$ws.Range("L2").Value = 622.1429000000001
$ws.Range("N2").Value = -848.1429000000001
$ws.Range("H43").Value = 3881.818
$ws.Range("I43").Value = 4500
$ws.Range("J43").Value = 3140
$ws.Range("K43").Value = 4500
$ws.Range("L43").Value = 3140
$ws.Range("M43").Value = -4431
$ws.Range("N43").Value = -3278
$ws.Range("H101").Value = 1167.2727
$ws.Range("I101").Value = 1263.2858
$ws.Range("K101").Value = 3789.8574
$ws.Range("M101").Value = -2167.8574
$ws.Range("H106").Value = 1212.8667
$ws.Range("I106").Value = 1245.9286
$ws.Range("K106").Value = 1245.9286
$ws.Range("M106").Value = -614.9286
$ws.Range("H137").Value = 3021.85
$ws.Range("I137").Value = 5104.4287
$ws.Range("J137").Value = 1900.4615
$ws.Range("K137").Value = 15313.2861
$ws.Range("L137").Value = 5701.3845
$ws.Range("M137").Value = -12763.2861
$ws.Range("N137").Value = -10801.3845
$ws.Range("H138").Value = 6806896
$ws.Range("J138").Value = 10758637
$ws.Range("L138").Value = 32275911
$ws.Range("N138").Value = -32286191

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 77975
$ws.Range("J80").Value = 77975
$ws.Range("L80").Value = 77975
$ws.Range("N80").Value = -79971
$ws.Range("H83").Value = 77975
$ws.Range("J83").Value = 77975
$ws.Range("L83").Value = 233925
$ws.Range("N83").Value = -243909
$ws.Range("H97").Value = 968.9697
$ws.Range("I97").Value = 1090.619
$ws.Range("J97").Value = 756.0833
$ws.Range("K97").Value = 1090.619
$ws.Range("L97").Value = 756.0833
$ws.Range("M97").Value = -594.6189999999999
$ws.Range("N97").Value = -1748.0833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2012.3
$ws.Range("I94").Value = 288
$ws.Range("K94").Value = 288
$ws.Range("M94").Value = 163
$ws.Range("H107").Value = 4971.3335
$ws.Range("I107").Value = 4967.875
$ws.Range("K107").Value = 4967.875
$ws.Range("M107").Value = -3047.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 276.5
$ws.Range("I11").Value = 150
$ws.Range("J11").Value = 403
$ws.Range("K11").Value = 150
$ws.Range("L11").Value = 403
$ws.Range("M11").Value = -10
$ws.Range("N11").Value = -683
$ws.Range("H31").Value = 4107.8438
$ws.Range("I31").Value = 3304.6428
$ws.Range("J31").Value = 4732.5557
$ws.Range("K31").Value = 3304.6428
$ws.Range("L31").Value = 4732.5557
$ws.Range("M31").Value = -3009.6428
$ws.Range("N31").Value = -5322.5557
$ws.Range("H34").Value = 4107.8438
$ws.Range("I34").Value = 3304.6428
$ws.Range("J34").Value = 4732.5557
$ws.Range("K34").Value = 3304.6428
$ws.Range("L34").Value = 4732.5557
$ws.Range("M34").Value = -3102.6428
$ws.Range("N34").Value = -5136.5557
$ws.Range("H105").Value = 9567.538
$ws.Range("I105").Value = 1232.8889
$ws.Range("K105").Value = 1232.8889
$ws.Range("M105").Value = 514.1111000000001
$ws.Range("H110").Value = 72993.5
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 72993.5
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 72993.5
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -81173.5
$ws.Range("H131").Value = 56058.6
$ws.Range("J131").Value = 56058.6
$ws.Range("L131").Value = 56058.6
$ws.Range("N131").Value = -66138.60000000001
$ws.Range("H134").Value = 2422.625
$ws.Range("I134").Value = 2212.4614
$ws.Range("K134").Value = 6637.3842
$ws.Range("M134").Value = -4102.3842
$ws.Range("H141").Value = 187846.28
$ws.Range("J141").Value = 187846.28
$ws.Range("L141").Value = 187846.28
$ws.Range("N141").Value = -198206.28

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 36490.816
$ws.Range("I56").Value = 36490.816
$ws.Range("K56").Value = 36490.816
$ws.Range("M56").Value = -35960.816
$ws.Range("H59").Value = 2380
$ws.Range("I59").Value = 2350
$ws.Range("K59").Value = 7050
$ws.Range("M59").Value = -6510
$ws.Range("H113").Value = 3360.875
$ws.Range("I113").Value = 2401.25
$ws.Range("J113").Value = 3680.75
$ws.Range("K113").Value = 7203.75
$ws.Range("L113").Value = 11042.25
$ws.Range("M113").Value = -5033.75
$ws.Range("N113").Value = -15382.25
$ws.Range("H131").Value = 1664.7142
$ws.Range("J131").Value = 1653.3823
$ws.Range("L131").Value = 4960.1469
$ws.Range("N131").Value = -15040.1469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2865.2666
$ws.Range("I80").Value = 2545
$ws.Range("J80").Value = 3345.6667
$ws.Range("K80").Value = 2545
$ws.Range("L80").Value = 3345.6667
$ws.Range("M80").Value = -1547
$ws.Range("N80").Value = -5341.6667
$ws.Range("H83").Value = 2865.2666
$ws.Range("I83").Value = 2545
$ws.Range("J83").Value = 3345.6667
$ws.Range("K83").Value = 12725
$ws.Range("L83").Value = 16728.3335
$ws.Range("M83").Value = -7733
$ws.Range("N83").Value = -26712.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 57497
$ws.Range("J74").Value = 58999.5
$ws.Range("L74").Value = 58999.5
$ws.Range("N74").Value = -60995.5
$ws.Range("H77").Value = 57497
$ws.Range("J77").Value = 58999.5
$ws.Range("L77").Value = 176998.5
$ws.Range("N77").Value = -186982.5
$ws.Range("H93").Value = 2662.6
$ws.Range("I93").Value = 1439
$ws.Range("K93").Value = 1439
$ws.Range("M93").Value = -191
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 74500
$ws.Range("J70").Value = 74500
$ws.Range("L70").Value = 74500
$ws.Range("N70").Value = -75130
$ws.Range("H73").Value = 74500
$ws.Range("J73").Value = 74500
$ws.Range("L73").Value = 74500
$ws.Range("N73").Value = -76684
$ws.Range("H81").Value = 1524.5333
$ws.Range("I81").Value = 1404.2
$ws.Range("K81").Value = 2808.4
$ws.Range("M81").Value = -1747.4
$ws.Range("H84").Value = 1524.5333
$ws.Range("I84").Value = 1404.2
$ws.Range("K84").Value = 14042
$ws.Range("M84").Value = -8738
$ws.Range("H122").Value = 65156.875
$ws.Range("I122").Value = 79041
$ws.Range("K122").Value = 237123
$ws.Range("M122").Value = -234673
